$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 108583.8495670803
$ws.Range("C2").Value = 205066.9665438474
$ws.Range("D2").Value = 277911.0873348299
$ws.Range("E2").Value = 308393.1483967978

$ws.Range("B3").Value = 132876.1585994226
$ws.Range("C3").Value = 246370.0193413306
$ws.Range("D3").Value = 328452.0744311343
$ws.Range("E3").Value = 366456.3985561453

$ws.Range("B4").Value = 133634.9505275033
$ws.Range("C4").Value = 247400.1052496259
$ws.Range("D4").Value = 328278.8617218114
$ws.Range("E4").Value = 365745.4199824685

$ws.Range("B5").Value = 106739.5658413006
$ws.Range("C5").Value = 183311.0458486494
$ws.Range("D5").Value = 231308.2449269884
$ws.Range("E5").Value = 257445.9193610857

$ws.Range("B6").Value = 92647.04710965729
$ws.Range("C6").Value = 160531.8596209518
$ws.Range("D6").Value = 204808.9336099969
$ws.Range("E6").Value = 225230.0381233434

$ws.Range("B7").Value = 10709.48881235689
$ws.Range("C7").Value = 17651.95417435372
$ws.Range("D7").Value = 22079.74300045208
$ws.Range("E7").Value = 23931.72673235543

$ws.Range("B8").Value = 235112.6018616934
$ws.Range("C8").Value = 595185.6476458332
$ws.Range("D8").Value = 917489.4878413631
$ws.Range("E8").Value = 1179526.881208779

$ws.Range("B9").Value = 124963.2271275988
$ws.Range("C9").Value = 220551.0157355429
$ws.Range("D9").Value = 287666.5132491607
$ws.Range("E9").Value = 320143.9981619639

$ws.Range("B10").Value = 59458.05984564151
$ws.Range("C10").Value = 99557.7727890665
$ws.Range("D10").Value = 127372.6658213741
$ws.Range("E10").Value = 136703.24863412

$ws.Range("B11").Value = 10986.64449520764
$ws.Range("C11").Value = 17410.58527303607
$ws.Range("D11").Value = 21997.24369316641
$ws.Range("E11").Value = 25471.83091001473

$ws.Range("B12").Value = 27332.21405171468
$ws.Range("C12").Value = 46896.137523741
$ws.Range("D12").Value = 58868.79312375189
$ws.Range("E12").Value = 62167.7280236506

$ws.Range("B13").Value = 33402.24755887037
$ws.Range("C13").Value = 55469.11180329873
$ws.Range("D13").Value = 71653.07652741244
$ws.Range("E13").Value = 77810.16062207166
